$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = 'Datos actualizados a 20 de Mayo de 2020 a las 19:05'

# Refreshed country statistics, re-sorted by Casos totales (column B) descending.
# Build the full target table (country + 7 numeric stats) and write it in one pass
# so the resulting row order and shared-string table exactly match the source update.
$data = @(
  ,@('Estados Unidos', 1576950, 6367, 364263, 1118829, 0, 325, 93858)
  ,@('Rusia', 308705, 8764, 85392, 220341, 0, 135, 2972)
  ,@('España', 278803, 0, 196958, 54067, 0, 0, 27778)
  ,@('Brasil', 275382, 3497, 106794, 150458, 0, 147, 18130)
  ,@('Reino Unido', 248818, 0, 0, 0, 0, 0, 35341)
  ,@('Italia', 226699, 0, 129401, 65129, 0, 0, 32169)
  ,@('Francia', 180809, 0, 62563, 90224, 0, 0, 28022)
  ,@('Alemania', 178170, 343, 156900, 13057, 0, 20, 8213)
  ,@('Turquia', 152587, 972, 113987, 34378, 0, 23, 4222)
  ,@('Iran', 126949, 2346, 98808, 20958, 0, 64, 7183)
  ,@('India', 111602, 5127, 45216, 62960, 0, 124, 3426)
  ,@('Peru', 99483, 0, 36524, 60045, 0, 0, 2914)
  ,@('China', 82965, 5, 78244, 87, 0, 0, 4634)
  ,@('Canada', 79503, 391, 40342, 33205, 0, 44, 5956)
  ,@('Arabia Saudita', 62545, 2691, 33478, 28728, 0, 10, 339)
  ,@('Belgica', 55983, 192, 14847, 31986, 0, 42, 9150)
  ,@('Mexico', 54346, 2713, 37325, 11355, 0, 334, 5666)
  ,@('Chile', 53617, 4038, 22504, 30569, 0, 35, 544)
  ,@('Pakistan', 45898, 1932, 13101, 31812, 0, 46, 985)
  ,@('Paises Bajos', 44447, 198, 0, 0, 0, 33, 5748)
  ,@('Catar', 37097, 1491, 6600, 30481, 0, 1, 16)
  ,@('Ecuador', 34151, 0, 3457, 27855, 0, 0, 2839)
  ,@('Bielorrusia', 32426, 918, 11415, 20832, 0, 4, 179)
  ,@('Suecia', 31523, 724, 4971, 22721, 0, 88, 3831)
  ,@('Suiza', 30658, 40, 27800, 966, 0, 1, 1892)
  ,@('Portugal', 29660, 228, 6452, 21945, 0, 16, 1263)
  ,@('Singapur', 29364, 570, 11207, 18135, 0, 0, 22)
  ,@('Banglades', 26738, 1617, 5207, 21145, 0, 16, 386)
  ,@('Emiratos Arabes Unidos', 25063, 0, 10791, 14045, 0, 0, 227)
  ,@('Irlanda', 24251, 0, 19470, 3220, 0, 0, 1561)
  ,@('Polonia', 19739, 471, 8183, 10594, 0, 14, 962)
  ,@('Ucrania', 19230, 354, 5955, 12711, 0, 16, 564)
  ,@('Indonesia', 19189, 693, 4575, 13372, 0, 21, 1242)
  ,@('Kuwait', 17568, 804, 4885, 12559, 0, 3, 124)
  ,@('Rumania', 17387, 196, 10356, 5884, 0, 10, 1147)
  ,@('Sudafrica', 17200, 0, 7960, 8928, 0, 0, 312)
  ,@('Colombia', 16935, 0, 4050, 12272, 0, 0, 613)
  ,@('Israel', 16667, 8, 13504, 2884, 0, 1, 279)
  ,@('Japon', 16367, 0, 11564, 4035, 0, 0, 768)
  ,@('Austria', 16353, 32, 14882, 838, 0, 1, 633)
  ,@('Egipto', 14229, 745, 3994, 9555, 0, 21, 680)
  ,@('Republica Dominicana', 13477, 254, 7142, 5889, 0, 5, 446)
  ,@('Filipinas', 13221, 279, 2932, 9447, 0, 5, 842)
  ,@('Dinamarca', 11117, 73, 9536, 1027, 0, 3, 554)
  ,@('Corea del Sur', 11110, 32, 10066, 781, 0, 0, 263)
  ,@('Serbia', 10833, 100, 5067, 5531, 0, 1, 235)
  ,@('Panama', 9867, 0, 6194, 3392, 0, 0, 281)
  ,@('Argentina', 8809, 0, 2933, 5482, 0, 1, 394)
  ,@('Chequia', 8698, 51, 5818, 2576, 0, 2, 304)
  ,@('Noruega', 8267, 0, 32, 8001, 0, 1, 234)
  ,@('Afganistan', 8145, 492, 930, 7028, 0, 9, 187)
  ,@('Barein', 7843, 311, 3334, 4497, 0, 0, 12)
  ,@('Argelia', 7542, 165, 3968, 3006, 0, 7, 568)
  ,@('Australia', 7079, 11, 6444, 535, 0, 0, 100)
  ,@('Marruecos', 7048, 25, 4037, 2817, 0, 1, 194)
  ,@('Malasia', 7009, 31, 5706, 1189, 0, 0, 114)
  ,@('Kazajistan', 6969, 218, 3649, 3285, 0, 0, 35)
  ,@('Moldavia', 6553, 213, 2953, 3372, 0, 7, 228)
  ,@('Finlandia', 6443, 44, 5000, 1139, 0, 3, 304)
  ,@('Nigeria', 6401, 0, 1734, 4475, 0, 0, 192)
  ,@('Ghana', 6096, 0, 1773, 4292, 0, 0, 31)
  ,@('Oman', 6043, 372, 1661, 4353, 0, 2, 29)
  ,@('Armenia', 5271, 230, 2419, 2785, 0, 3, 67)
  ,@('Bolivia', 4481, 218, 533, 3759, 0, 15, 189)
  ,@('Luxemburgo', 3958, 0, 3718, 131, 0, 0, 109)
  ,@('Irak', 3724, 113, 2438, 1152, 0, 3, 134)
  ,@('Azerbaiyan', 3631, 113, 2253, 1335, 0, 2, 43)
  ,@('Hungria', 3598, 42, 1454, 1674, 0, 3, 470)
  ,@('Camerun', 3529, 0, 1567, 1822, 0, 0, 140)
  ,@('Tailandia', 3034, 1, 2888, 90, 0, 0, 56)
  ,@('Honduras', 2955, 157, 349, 2459, 0, 1, 147)
  ,@('Uzbekistan', 2927, 72, 2369, 545, 0, 0, 13)
  ,@('Guinea', 2863, 0, 1525, 1320, 0, 0, 18)
  ,@('Grecia', 2850, 10, 1374, 1310, 0, 1, 166)
  ,@('Sudan', 2728, 0, 286, 2331, 0, 0, 111)
  ,@('Senegal', 2714, 97, 1186, 1498, 0, 0, 30)
  ,@('Bosnia y Herzegovina', 2338, 17, 1557, 645, 0, 2, 136)
  ,@('Bulgaria', 2292, 33, 684, 1492, 0, 4, 116)
  ,@('Croacia', 2234, 2, 1978, 160, 0, 0, 96)
  ,@('Costa de Marfil', 2153, 0, 1050, 1075, 0, 0, 28)
  ,@('Tayikistan', 2140, 204, 641, 1458, 0, 0, 41)
  ,@('Guatemala', 2133, 132, 155, 1935, 0, 5, 43)
  ,@('Cuba', 1900, 13, 1573, 248, 0, 0, 79)
  ,@('Republica de Macedonia', 1858, 19, 1367, 381, 0, 4, 110)
  ,@('Republica de Yibuti', 1828, 210, 1052, 767, 0, 2, 9)
  ,@('Islandia', 1803, 1, 1789, 4, 0, 0, 10)
  ,@('Estonia', 1794, 3, 956, 774, 0, 0, 64)
  ,@('Consejo Danes para los Refugiados', 1731, 102, 302, 1368, 0, 0, 61)
  ,@('Lituania', 1577, 15, 1049, 468, 0, 0, 60)
  ,@('Somalia', 1573, 71, 188, 1324, 0, 2, 61)
  ,@('El Salvador', 1571, 73, 531, 1009, 0, 1, 31)
  ,@('Nueva Zelanda', 1503, 0, 1447, 35, 0, 0, 21)
  ,@('Gabon', 1502, 0, 318, 1172, 0, 0, 12)
  ,@('Eslovaquia', 1496, 1, 1231, 237, 0, 0, 28)
  ,@('Mayotte', 1475, 105, 894, 562, 0, 1, 19)
  ,@('Eslovenia', 1468, 1, 1340, 23, 0, 1, 105)
  ,@('Kirguistan', 1270, 27, 910, 346, 0, 0, 14)
  ,@('Maldivas', 1186, 43, 91, 1091, 0, 0, 4)
  ,@('Hong Kong', 1056, 0, 1026, 26, 0, 0, 4)
  ,@('Tunez', 1044, 1, 826, 171, 0, 1, 47)
  ,@('Guinea-Bisau', 1038, 0, 42, 990, 0, 0, 6)
  ,@('Kenia', 1029, 66, 366, 613, 0, 0, 50)
  ,@('Sri Lanka', 1027, 4, 584, 434, 0, 0, 9)
  ,@('Letonia', 1016, 4, 694, 301, 0, 0, 21)
  ,@('Albania', 964, 15, 758, 175, 0, 0, 31)
  ,@('Libano', 961, 7, 251, 684, 0, 0, 26)
  ,@('Republica de Chipre', 922, 4, 516, 389, 0, 0, 17)
  ,@('Niger', 914, 0, 734, 125, 0, 0, 55)
  ,@('Mali', 901, 0, 529, 319, 0, 0, 53)
  ,@('Costa Rica', 882, 0, 577, 295, 0, 0, 10)
  ,@('Zambia', 832, 60, 197, 628, 0, 0, 7)
  ,@('Paraguay', 829, 0, 230, 588, 0, 0, 11)
  ,@('Guinea Ecuatorial', 825, 0, 22, 796, 0, 0, 7)
  ,@('Burkina Faso', 809, 13, 661, 96, 0, 1, 52)
  ,@('Principado de Andorra', 761, 0, 628, 82, 0, 0, 51)
  ,@('Venezuela', 749, 131, 253, 486, 0, 0, 10)
  ,@('Uruguay', 738, 0, 579, 139, 0, 0, 20)
  ,@('Georgia', 713, 6, 475, 226, 0, 0, 12)
  ,@('Crucero', 712, 0, 651, 48, 0, 0, 13)
  ,@('Jordania', 672, 23, 446, 217, 0, 0, 9)
  ,@('San Marino', 656, 1, 220, 395, 0, 0, 41)
  ,@('Haiti', 596, 63, 21, 553, 0, 1, 22)
  ,@('Malta', 584, 15, 465, 113, 0, 0, 6)
  ,@('Sierra Leona', 570, 36, 205, 331, 0, 1, 34)
  ,@('Republica del Chad', 545, 0, 139, 350, 0, 0, 56)
  ,@('Jamaica', 520, 0, 145, 366, 0, 0, 9)
  ,@('Tanzania', 509, 0, 183, 305, 0, 0, 21)
  ,@('Reunion', 447, 1, 411, 35, 0, 1, 1)
  ,@('Taiwan', 440, 0, 402, 31, 0, 0, 7)
  ,@('Nepal', 427, 25, 45, 380, 0, 0, 2)
  ,@('Congo', 420, 0, 132, 273, 0, 0, 15)
  ,@('Estado de Palestina', 398, 7, 346, 50, 0, 0, 2)
  ,@('Etiopia', 389, 24, 122, 262, 0, 0, 5)
  ,@('Madagascar', 371, 45, 131, 238, 0, 0, 2)
  ,@('Republica de Africa Central', 366, 0, 18, 348, 0, 0, 0)
  ,@('Cabo Verde', 349, 14, 85, 261, 0, 0, 3)
  ,@('Togo', 338, 0, 107, 219, 0, 0, 12)
  ,@('Isla de Man', 336, 1, 302, 10, 0, 0, 24)
  ,@('Mauricio', 332, 0, 322, 0, 0, 0, 10)
  ,@('Vietnam', 324, 0, 263, 61, 0, 0, 0)
  ,@('Montenegro', 324, 0, 312, 3, 0, 0, 9)
  ,@('Ruanda', 308, 0, 209, 99, 0, 0, 0)
  ,@('Sudan del Sur', 290, 0, 4, 282, 0, 0, 4)
  ,@('Uganda', 264, 4, 65, 199, 0, 0, 0)
  ,@('Nicaragua', 254, 0, 199, 38, 0, 0, 17)
  ,@('Santo Tome y Principe', 251, 0, 4, 239, 0, 0, 8)
  ,@('Liberia', 238, 5, 128, 87, 0, 0, 23)
  ,@('Guayana Francesa', 218, 0, 136, 81, 0, 0, 1)
  ,@('Suazilandia', 208, 0, 87, 119, 0, 0, 2)
  ,@('Birmania', 193, 0, 104, 83, 0, 0, 6)
  ,@('Martinica', 192, 0, 91, 87, 0, 0, 14)
  ,@('Islas Feroe', 187, 0, 187, 0, 0, 0, 0)
  ,@('Yemen', 167, 0, 5, 134, 0, 0, 28)
  ,@('Mozambique', 156, 10, 48, 108, 0, 0, 0)
  ,@('Guadalupe', 155, 0, 109, 33, 0, 0, 13)
  ,@('Gibraltar', 149, 2, 145, 4, 0, 0, 0)
  ,@('Brunei', 141, 0, 136, 4, 0, 0, 1)
  ,@('Mongolia', 140, 0, 26, 114, 0, 0, 0)
  ,@('Mauritania', 131, 0, 7, 120, 0, 0, 4)
  ,@('Benin', 130, 0, 83, 45, 0, 0, 2)
  ,@('Guyana', 125, 0, 46, 69, 0, 0, 10)
  ,@('Bermudas', 125, 0, 78, 38, 0, 0, 9)
  ,@('Camboya', 122, 0, 122, 0, 0, 0, 0)
  ,@('Trinidad yTobago', 116, 0, 107, 1, 0, 0, 8)
  ,@('Islas Caimanes', 111, 0, 55, 55, 0, 0, 1)
  ,@('Aruba', 101, 0, 93, 5, 0, 0, 3)
  ,@('Monaco', 97, 0, 87, 6, 0, 0, 4)
  ,@('Bahamas', 96, 0, 43, 42, 0, 0, 11)
  ,@('Barbados', 90, 2, 68, 15, 0, 0, 7)
  ,@('Liechtenstein', 82, 0, 55, 26, 0, 0, 1)
  ,@('San Martin (Parte Holandesa)', 77, 0, 54, 8, 0, 0, 15)
  ,@('Malaui', 71, 1, 27, 41, 0, 0, 3)
  ,@('Libia', 69, 1, 35, 31, 0, 0, 3)
  ,@('Polinesia Francesa', 60, 0, 60, 0, 0, 0, 0)
  ,@('Siria', 58, 0, 36, 19, 0, 0, 3)
  ,@('Angola', 52, 0, 17, 32, 0, 0, 3)
  ,@('Zimbabue', 46, 0, 18, 24, 0, 0, 4)
  ,@('Macao', 45, 0, 45, 0, 0, 0, 0)
  ,@('Burundi', 42, 0, 20, 21, 0, 0, 1)
  ,@('San Martin (Parte Francesa)', 40, 1, 33, 4, 0, 0, 3)
  ,@('Puerto Rico', 39, 0, 1, 36, 0, 0, 2)
  ,@('Eritrea', 39, 0, 39, 0, 0, 0, 0)
  ,@('Comoras', 34, 23, 8, 25, 0, 0, 1)
  ,@('Guam', 32, 0, 0, 31, 0, 0, 1)
  ,@('Botsuana', 25, 0, 17, 7, 0, 0, 1)
  ,@('Antigua y Barbuda', 25, 0, 19, 3, 0, 0, 3)
  ,@('Gambia', 24, 0, 13, 10, 0, 0, 1)
  ,@('Timor Oriental', 24, 0, 24, 0, 0, 0, 0)
  ,@('Granada', 22, 0, 14, 8, 0, 0, 0)
  ,@('Butan', 21, 0, 5, 16, 0, 0, 0)
  ,@('Laos', 19, 0, 14, 5, 0, 0, 0)
  ,@('Fiyi', 18, 0, 15, 3, 0, 0, 0)
  ,@('Santa Lucia', 18, 0, 18, 0, 0, 0, 0)
  ,@('Belice', 18, 0, 16, 0, 0, 0, 2)
  ,@('Nueva Caledonia', 18, 0, 18, 0, 0, 0, 0)
  ,@('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0)
  ,@('San Vicente y las Granadinas', 17, 0, 14, 3, 0, 0, 0)
  ,@('Namibia', 16, 0, 14, 2, 0, 0, 0)
  ,@('Curazao', 16, 0, 14, 1, 0, 0, 1)
  ,@('Dominica', 16, 0, 16, 0, 0, 0, 0)
  ,@('San Cristobal y Nieves', 15, 0, 15, 0, 0, 0, 0)
  ,@('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0)
  ,@('Santa Sede', 12, 0, 2, 10, 0, 0, 0)
  ,@('Islas Turcas y Caicos', 12, 0, 10, 1, 0, 0, 1)
  ,@('Surinam', 11, 0, 9, 1, 0, 0, 1)
  ,@('Montserrat', 11, 0, 10, 0, 0, 0, 1)
  ,@('Groenlandia', 11, 0, 11, 0, 0, 0, 0)
  ,@('Seychelles', 11, 0, 11, 0, 0, 0, 0)
  ,@('Islas Virgenes Britanicas', 8, 0, 6, 1, 0, 0, 1)
  ,@('Papua Nueva Guinea', 8, 0, 8, 0, 0, 0, 0)
  ,@('Bonaire, San Eustaquio y Saba', 6, 0, 6, 0, 0, 0, 0)
  ,@('Sahara Occidental', 6, 0, 6, 0, 0, 0, 0)
  ,@('San Bartolome', 6, 0, 6, 0, 0, 0, 0)
  ,@('Anguila', 3, 0, 3, 0, 0, 0, 0)
  ,@('Lesoto', 1, 0, 0, 1, 0, 0, 0)
  ,@('San Pedro y Miquelon', 1, 0, 1, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 4 + $i
    $row = $data[$i]
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Output "Updated $($data.Length) country rows"
